$d = $word.ActiveDocument
$laquo = [char]0x00AB
$raquo = [char]0x00BB
$nbsp  = [char]0x00A0
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- Paragraph 5: "HTML Comments<nbsp>: " -> "HTML " + "Comments" (spell-wrapped) + "<nbsp>: " ---
$p5 = $d.Paragraphs(5)
$r5 = $p5.Range
$xml5 = '<w:p ' + $wns + '><w:r><w:t xml:space="preserve">HTML </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Comments</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">' + $nbsp + ': </w:t></w:r></w:p>'
$r5.InsertXML($xml5)

# --- Paragraph 1: "Docx Comments<nbsp>: " -> new paragraph "Inline HTML<nbsp>: " + mergefield, then "Docx Comments<nbsp>: " ---
$p1 = $d.Paragraphs(1)
$r1 = $p1.Range
$xml1 = '<w:p ' + $wns + '><w:proofErr w:type="spellStart"/><w:r><w:t>Inline</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> HTML' + $nbsp + ': </w:t></w:r><w:r><w:fldChar w:fldCharType="begin"/></w:r><w:r><w:instrText xml:space="preserve"> MERGEFIELD  $comments_html  \* MERGEFORMAT </w:instrText></w:r><w:r><w:fldChar w:fldCharType="separate"/></w:r><w:r><w:rPr><w:noProof/></w:rPr><w:t>' + $laquo + '$comments_html' + $raquo + '</w:t></w:r><w:r><w:fldChar w:fldCharType="end"/></w:r></w:p><w:p ' + $wns + '><w:proofErr w:type="spellStart"/><w:r><w:t>Docx</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Comments</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">' + $nbsp + ': </w:t></w:r></w:p>'
$r1.InsertXML($xml1)

# --- Footer: "My Footer" -> "My" + " " + "Footer" (spell-wrapped), regular space (no colon, no NBSP) ---
$sec = $d.Sections(1)
$ftr = $sec.Footers(1)
$rf = $ftr.Range
$xmlf = '<w:p ' + $wns + '><w:pPr><w:pStyle w:val="Pieddepage"/><w:jc w:val="right"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>My</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Footer</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
$rf.InsertXML($xmlf)

Write-Output "done"
